$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.100.94"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.526.52"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.07%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.09"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +1.16%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.02"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "3.525.34"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("E8").Value = "  +0.07%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  +2.03%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.78"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -7.56%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("D13").Value = "4.125.98"
$ws.Range("E13").Value = "  +2.66%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000184"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +1.67%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.01"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "3.530.21"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "65.209.60"
$ws.Range("E18").Value = "  -0.49%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.28"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +4.22%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.56"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.31%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.571"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").Value = "3.668.39"
$ws.Range("E24").Value = "  +2.68%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.74"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.49%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +0.00%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +6.99%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.66"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +6.41%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +2.12%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.15"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "3.540.46"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("E33").Value = "  +0.01%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.73"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("E35").Value = "  -0.52%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.25"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +7.09%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.92"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +0.32%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.40"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  +4.21%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.96"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +3.04%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0796"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +4.30%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.822"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.36%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.92"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +13.96%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.76"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.32%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +0.10%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.41"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -0.42%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.66"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +2.52%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +4.46%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.77"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").Value = "2.379.08"
$ws.Range("E50").Value = "  +7.57%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "301.04"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +6.18%  "
